# Auto-generated edit script: updates cryptocurrency price/volume data
# per the GitHub Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.859.06'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.894.69'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '0.7922'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').Value = '242.59'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '0.3204'
$ws.Range('E8').Value = '  +2.55%  '
$ws.Range('D9').Value = '26.19'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '0.07104'
$ws.Range('E10').Value = '  +3.32%  '
$ws.Range('D11').Value = '0.08063'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '0.7736'
$ws.Range('E12').Value = '  +5.11%  '
$ws.Range('D13').Value = '1.913.84'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.320'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('D15').Value = '92.41'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '29.902.94'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '13.89'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '5.933'
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').Value = '244.51'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '0.000007764'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = '2.165.32'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = '8.207'
$ws.Range('E22').Value = '  +19.94%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '0.1624'
$ws.Range('E25').Value = '  +14.16%  '
$ws.Range('D26').Value = '9.329'
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('D27').Value = '164.63'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').Value = '18.71'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = '2.077'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('D30').Value = '1.377'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').Value = '1.541'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('D32').Value = '4.486'
$ws.Range('E32').Value = '  +4.80%  '
$ws.Range('D33').Value = '0.05661'
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.110'
$ws.Range('E34').Value = '  +1.24%  '
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('D36').Value = '0.7396'
$ws.Range('E36').Value = '  +2.10%  '
$ws.Range('D37').Value = '1.003'
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').Value = '2.697'
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = '2.782'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').Value = '0.4461'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').Value = '72.36'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').Value = '5.875'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '0.8469'
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('D45').Value = '1.002'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '1.893'
$ws.Range('E46').Value = '  +2.28%  '
$ws.Range('D47').Value = '1.030.27'
$ws.Range('E47').Value = '  +5.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.50'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('D49').Value = '9.911'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('D50').Value = '7.503'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.059.66'
$ws.Range('E51').Value = '  -0.05%  '
